$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Merge the five separate runs that make up the G5 paragraph text
#    into a single run, by replacing the whole paragraph text with
#    itself via Find/Replace (Word recreates it as one run).
# ------------------------------------------------------------------
$g5Text = "G5: Authorities should have access to refined data related to committed violation"
$d.Content.Find.Execute($g5Text, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $g5Text, 2)

# ------------------------------------------------------------------
# 2) Remove the existing "_GoBack" bookmark (it currently sits at the
#    end of the G5 paragraph) - it will be recreated at the end of the
#    new G6 paragraph below.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 3) Find the G5 paragraph and insert a brand-new paragraph after it
#    containing the G6 text. A placeholder character "X" is appended
#    temporarily so that the position right after the real G6 text is
#    NOT the paragraph-end (pilcrow) position while we add the
#    bookmark there (adding a collapsed bookmark exactly at a
#    paragraph-end position is mishandled by this host and snaps back
#    to the start of the document).
# ------------------------------------------------------------------
$g6Text = "G6: Authorities should be able to communicate accident reports"

$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq ($g5Text + "`r")) {
        $rng = $p.Range
        $rng.Collapse(0)
        $rng.InsertParagraphAfter()
        $rng.InsertAfter($g6Text + "X")

        # Collapsed range sitting right after the real G6 text, before
        # the temporary "X" placeholder (i.e. not at the pilcrow).
        $bmRange = $rng.Duplicate
        $bmRange.MoveEnd(1, -1)
        $bmRange.Collapse(0)
        $d.Bookmarks.Add("_GoBack", $bmRange)

        # Remove the temporary placeholder character now that the
        # bookmark has been safely anchored: start from the bookmark's
        # own (collapsed) position and extend one character forward to
        # cover the "X" placeholder, then delete it.
        $padRange = $d.Bookmarks.Item("_GoBack").Range.Duplicate
        $padRange.Collapse(0)
        $padRange.MoveEnd(1, 1)
        $padRange.Delete()
        break
    }
}
